# Refresh the cryptos snapshot on Sheet1 with the latest scraped values.
#
# Column D entries that look like plain numbers (e.g. "209.41") are forced
# to stay literal text -- exactly like the rest of the sheet -- by building
# the string with a leading apostrophe ("'" + "209.41"), the same trick used
# when typing such a value straight into Excel so it is not reinterpreted as
# a number. Thousands-grouped prices (e.g. "25.842.29") are already unambiguous
# text and do not need the prefix.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.842.29"
$ws.Range("E2").Value = "  -0.48%  "

$ws.Range("D3").Value = "1.601.28"
$ws.Range("E3").Value = "  -2.10%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "'" + "209.41"
$ws.Range("E5").Value = "  -2.36%  "

$ws.Range("E6").Value = "  -0.15%  "

$ws.Range("D7").Value = "'" + "0.480"
$ws.Range("E7").Value = "  -4.80%  "

$ws.Range("D8").Value = "'" + "0.247"
$ws.Range("E8").Value = "  -2.28%  "

$ws.Range("D9").Value = "'" + "0.0610"
$ws.Range("E9").Value = "  -2.22%  "

$ws.Range("D10").Value = "'" + "17.94"
$ws.Range("E10").Value = "  -3.19%  "

$ws.Range("D11").Value = "'" + "0.0784"
$ws.Range("E11").Value = "  -0.88%  "

$ws.Range("D12").Value = "1.822.64"
$ws.Range("E12").Value = "  -2.20%  "

$ws.Range("D13").Value = "1.602.52"
$ws.Range("E13").Value = "  -2.04%  "

$ws.Range("E14").Value = "  -3.05%  "

$ws.Range("D15").Value = "'" + "0.509"
$ws.Range("E15").Value = "  -4.46%  "

$ws.Range("D16").Value = "25.833.72"
$ws.Range("E16").Value = "  -0.58%  "

$ws.Range("D17").Value = "'" + "60.51"
$ws.Range("E17").Value = "  -2.01%  "

$ws.Range("E19").Value = "  -0.09%  "

$ws.Range("D20").Value = "'" + "189.53"
$ws.Range("E20").Value = "  -0.72%  "

$ws.Range("D21").Value = "'" + "4.18"
$ws.Range("E21").Value = "  -1.58%  "

$ws.Range("D22").Value = "'" + "9.33"
$ws.Range("E22").Value = "  -3.23%  "

$ws.Range("D23").Value = "'" + "5.95"
$ws.Range("E23").Value = "  -2.86%  "

$ws.Range("B24").Value = "Monero"
$ws.Range("C24").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D24").Value = "'" + "141.83"
$ws.Range("E24").Value = "  -1.16%  "

$ws.Range("B25").Value = "BinanceUSD"
$ws.Range("C25").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D25").Value = "'" + "1.01"
$ws.Range("E25").Value = "  -0.10%  "

$ws.Range("D26").Value = "'" + "0.127"
$ws.Range("E26").Value = "  -4.61%  "

$ws.Range("E27").Value = "  -3.80%  "

$ws.Range("D28").Value = "'" + "6.51"

$ws.Range("D29").Value = "'" + "14.95"
$ws.Range("E29").Value = "  -2.06%  "

$ws.Range("D30").Value = "'" + "1.19"
$ws.Range("E30").Value = "  -3.94%  "

$ws.Range("D31").Value = "'" + "0.0467"
$ws.Range("E31").Value = "  -3.43%  "

$ws.Range("D32").Value = "'" + "3.08"
$ws.Range("E32").Value = "  -2.47%  "

$ws.Range("D33").Value = "'" + "3.00"
$ws.Range("E33").Value = "  -4.88%  "

$ws.Range("E34").Value = "  -0.68%  "

$ws.Range("D35").Value = "'" + "1.47"
$ws.Range("E35").Value = "  -1.75%  "

$ws.Range("D36").Value = "1.098.79"
$ws.Range("E36").Value = "  -3.42%  "

$ws.Range("D38").Value = "'" + "1.00"
$ws.Range("E38").Value = "  -0.19%  "

$ws.Range("D39").Value = "'" + "0.793"
$ws.Range("E39").Value = "  -8.50%  "

$ws.Range("E40").Value = "  -2.56%  "

$ws.Range("D41").Value = "'" + "0.497"
$ws.Range("E41").Value = "  -4.96%  "

$ws.Range("E42").Value = "  -2.94%  "

$ws.Range("D43").Value = "1.734.08"
$ws.Range("E43").Value = "  -2.21%  "

$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "'" + "0.743"
$ws.Range("E44").Value = "  -4.53%  "

$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'" + "5.05"
$ws.Range("E45").Value = "  -3.76%  "

$ws.Range("E46").Value = "  -1.82%  "

$ws.Range("D47").Value = "'" + "53.18"
$ws.Range("E47").Value = "  -3.89%  "

$ws.Range("E48").Value = "  -3.22%  "

$ws.Range("D49").Value = "'" + "1.44"
$ws.Range("E49").Value = "  -3.89%  "

$ws.Range("D50").Value = "'" + "0.409"
$ws.Range("E50").Value = "  -1.26%  "

$ws.Range("E51").Value = "  -0.24%  "
